$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 11: new timesheet entry ---
# Copy the date-formatted style from A10 onto A11 so it matches the other
# date cells (s="1"), then set the actual date value (2016-02-02 = 42402).
$ws.Range("A10").Copy() | Out-Null
$ws.Range("A11").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("A11").Value2 = 42402

$ws.Range("B11").Value2 = 1
$ws.Range("C11").Value2 = "Nieuws items editen (start)"

# --- E5: totaal range grows to include the new row ---
$ws.Range("E5").Formula = "=SUM(B4:B29)"

# --- Selection moved to C17 ---
$ws.Range("C17").Select() | Out-Null
